$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "-"

$ws.Range("B3").Value = "-"
$ws.Range("D3").Value = "-"

$ws.Range("C4").Value = "MEC-3B-Retificação"
$ws.Range("E4").Value = "MCT-3A-Processos de Usinagem 2"

$ws.Range("C6").Value = "MEC-3B-Retificação"
$ws.Range("E6").Value = "MCT-3A-Processos de Usinagem 2"

$ws.Range("C7").Value = "MEC-3B-Retificação"
$ws.Range("E7").Value = "MCT-3A-Processos de Usinagem 2"

$ws.Range("C8").Value = "MEC-3B-Retificação"
$ws.Range("E8").Value = "MCT-3A-Processos de Usinagem 2"

$ws.Range("D14").Value = "MEC-3A-Elemaq."
$ws.Range("F14").Value = "-"

$ws.Range("D15").Value = "MEC-3A-Elemaq."
$ws.Range("F15").Value = "-"
